$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HISAT2")

# --- Fill in the missing values for row 7 (SRP179837_leaf3) ---
$ws.Range("B7").Value = 38809128
$ws.Range("C7").Value = 0.9191
$ws.Range("C7").NumberFormat = "0.00%"
$ws.Range("D7").Value = 0.7769
$ws.Range("D7").NumberFormat = "0.00%"
$ws.Range("E7").Value = 0.106
$ws.Range("E7").NumberFormat = "0.00%"

# --- Add a new "With trimmomatic applied first" table below the existing one ---
$ws.Range("A18").Value = "With trimmomatic applied first to remove adapter sequences"
$ws.Range("A19").Value = "Parameters set: none"

$ws.Range("A20").Value = "Sample"
$ws.Range("B20").Value = "# of input reads"
$ws.Range("C20").Value = "Overall alignment rate"
$ws.Range("D20").Value = "% uniquely mapped reads"
$ws.Range("E20").Value = "% multimapping reads"
$ws.Range("A20:E20").Font.Bold = $true

$ws.Range("A21").Value = "SRP179837_leaf1"
$ws.Range("A22").Value = "SRP179837_leaf2"
$ws.Range("A23").Value = "SRP179837_leaf3"
$ws.Range("A24").Value = "SRP179837_root1"
$ws.Range("A25").Value = "SRP179837_root2"
$ws.Range("A26").Value = "SRP179837_root3"
$ws.Range("A27").Value = "RAC_leaf1"
$ws.Range("A28").Value = "RAC_leaf2"
$ws.Range("A29").Value = "RAC_leaf3"
$ws.Range("A30").Value = "RAC_root1"
$ws.Range("A31").Value = "RAC_root2"
$ws.Range("A32").Value = "RAC_root3"

# --- Update the selection to match the final cursor position ---
$ws.Range("F20").Select() | Out-Null
